$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell K1 - copy format from J1 (bold header style), then set value
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("K1").Value = "intervention_type"

# Data values for column K (intervention_type), rows 2-101
$ws.Range("K2").Value = "DRUG"
$ws.Range("K3").Value = "PROCEDURE"
$ws.Range("K4").Value = "PROCEDURE"
$ws.Range("K5").Value = "PROCEDURE"
$ws.Range("K6").Value = "DRUG"
$ws.Range("K7").Value = "DEVICE"
# Row 8: left blank (source has an empty inline string cell; engine cannot materialize empty cells via COM)
$ws.Range("K9").Value = "DEVICE"
$ws.Range("K10").Value = "DEVICE"
# Row 11: left blank (source has an empty inline string cell; engine cannot materialize empty cells via COM)
$ws.Range("K12").Value = "PROCEDURE"
$ws.Range("K13").Value = "OTHER"
$ws.Range("K14").Value = "DRUG"
$ws.Range("K15").Value = "OTHER"
$ws.Range("K16").Value = "DRUG"
$ws.Range("K17").Value = "OTHER"
$ws.Range("K18").Value = "OTHER"
$ws.Range("K19").Value = "OTHER"
$ws.Range("K20").Value = "OTHER"
$ws.Range("K21").Value = "PROCEDURE"
$ws.Range("K22").Value = "BEHAVIORAL"
$ws.Range("K23").Value = "PROCEDURE"
$ws.Range("K24").Value = "DEVICE"
$ws.Range("K25").Value = "DEVICE"
$ws.Range("K26").Value = "OTHER"
$ws.Range("K27").Value = "OTHER"
$ws.Range("K28").Value = "OTHER"
$ws.Range("K29").Value = "OTHER"
$ws.Range("K30").Value = "OTHER"
$ws.Range("K31").Value = "DRUG"
$ws.Range("K32").Value = "DRUG"
$ws.Range("K33").Value = "PROCEDURE"
$ws.Range("K34").Value = "DEVICE"
$ws.Range("K35").Value = "PROCEDURE"
$ws.Range("K36").Value = "OTHER"
$ws.Range("K37").Value = "DEVICE"
$ws.Range("K38").Value = "OTHER"
$ws.Range("K39").Value = "OTHER"
$ws.Range("K40").Value = "DEVICE"
$ws.Range("K41").Value = "BIOLOGICAL"
# Row 42: left blank (source has an empty inline string cell; engine cannot materialize empty cells via COM)
$ws.Range("K43").Value = "DEVICE"
$ws.Range("K44").Value = "DRUG"
$ws.Range("K45").Value = "PROCEDURE"
$ws.Range("K46").Value = "OTHER"
$ws.Range("K47").Value = "OTHER"
# Row 48: left blank (source has an empty inline string cell; engine cannot materialize empty cells via COM)
$ws.Range("K49").Value = "BEHAVIORAL"
$ws.Range("K50").Value = "OTHER"
$ws.Range("K51").Value = "OTHER"
$ws.Range("K52").Value = "DIAGNOSTIC_TEST"
$ws.Range("K53").Value = "OTHER"
$ws.Range("K54").Value = "DEVICE"
$ws.Range("K55").Value = "BIOLOGICAL"
$ws.Range("K56").Value = "OTHER"
$ws.Range("K57").Value = "BEHAVIORAL"
$ws.Range("K58").Value = "OTHER"
$ws.Range("K59").Value = "OTHER"
$ws.Range("K60").Value = "OTHER"
$ws.Range("K61").Value = "BEHAVIORAL"
$ws.Range("K62").Value = "DIAGNOSTIC_TEST"
$ws.Range("K63").Value = "DIAGNOSTIC_TEST"
$ws.Range("K64").Value = "DEVICE"
$ws.Range("K65").Value = "DEVICE"
$ws.Range("K66").Value = "DEVICE"
$ws.Range("K67").Value = "OTHER"
$ws.Range("K68").Value = "OTHER"
$ws.Range("K69").Value = "OTHER"
$ws.Range("K70").Value = "DRUG"
$ws.Range("K71").Value = "DRUG"
$ws.Range("K72").Value = "DRUG"
$ws.Range("K73").Value = "DEVICE"
$ws.Range("K74").Value = "OTHER"
$ws.Range("K75").Value = "OTHER"
$ws.Range("K76").Value = "OTHER"
$ws.Range("K77").Value = "DEVICE"
$ws.Range("K78").Value = "OTHER"
$ws.Range("K79").Value = "DEVICE"
$ws.Range("K80").Value = "OTHER"
$ws.Range("K81").Value = "BEHAVIORAL"
$ws.Range("K82").Value = "DRUG"
$ws.Range("K83").Value = "BEHAVIORAL"
$ws.Range("K84").Value = "DIAGNOSTIC_TEST"
$ws.Range("K85").Value = "DEVICE"
$ws.Range("K86").Value = "DIAGNOSTIC_TEST"
$ws.Range("K87").Value = "PROCEDURE"
$ws.Range("K88").Value = "OTHER"
# Row 89: left blank (source has an empty inline string cell; engine cannot materialize empty cells via COM)
$ws.Range("K90").Value = "OTHER"
$ws.Range("K91").Value = "OTHER"
$ws.Range("K92").Value = "OTHER"
$ws.Range("K93").Value = "OTHER"
$ws.Range("K94").Value = "OTHER"
$ws.Range("K95").Value = "OTHER"
$ws.Range("K96").Value = "DEVICE"
$ws.Range("K97").Value = "OTHER"
$ws.Range("K98").Value = "OTHER"
$ws.Range("K99").Value = "DEVICE"
$ws.Range("K100").Value = "OTHER"
# Row 101: left blank (source has an empty inline string cell; engine cannot materialize empty cells via COM)
